# CONTACTOS-LA-INTEGRAL.xlsx — "Partner with new structure for upload"
#
# The authoritative content change in this commit is the text of the
# shared string used by B2 (row 2, "nombre" column): the placeholder
# "María Modificada" is replaced with "Sin modificar a Marías". The
# author's last recorded selection in the sheet was also B2 (previously
# B3), so we move the active cell to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the contact record's "nombre" value for row 2.
$ws.Range("B2").Value = "Sin modificar a Marías"

# Match the saved cursor/selection state (B2 instead of B3).
[void]$ws.Range("B2").Select()

# Workbook was authored with the (non-1904) date system; keep it explicit.
$wb.Date1904 = $false
